$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F2: "Kode Fixed Income : Hasil Generate" -> "Kode Fixed Income : OBL00108"
$ws.Range("F2").Value = "Username : 32070;`nPassword : bni1234;`nRole : 18/19 - Pimpinan Kelompok Investasi/Pengelola Investasi;`nKode Fixed Income : OBL00108"

# Update M2: "OBL00107" -> "OBL00108"
$ws.Range("M2").Value = "OBL00108"

# Row 2 height: 90 -> 75
$ws.Rows.Item(2).RowHeight = 75
